# BBI-23-1.xlsx — "Add files via upload" edit
#
# Updates a handful of grade-sheet input cells (columns E/G/H feed the
# MAX(...) formula in column J, which in turn drives K/L), fixes a cell
# that had erroneously been typed as text ("4 (без допуска)") back into a
# plain number, and adds a new (empty, bold-formatted) cell a few rows
# below the table — exactly as captured by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Строка 8 (Губеева София Юрьевна): Лаба №3 -2 -> 3
$ws.Range("E8").Value = 3

# Строка 10 (Ибрамхалилов Роман Ламетович): Лаба №5 -1 -> 4, КР-1 (пусто) -> 0
$ws.Range("G10").Value = 4
$ws.Range("H10").Value = 0

# Строка 11 (Капелина Анна Ивановна): Лаба №3 -2 -> 3
$ws.Range("E11").Value = 3

# Строка 20 (Стоценко Александр Сергеевич): КР-1 0 -> 5
$ws.Range("H20").Value = 5

# Строка 21 (Таранов Артём Игоревич): Лаба №3 was text "4 (без допуска)" -> numeric 4
$ws.Range("E21").Value = 4

# Строка 22 (Титова Надежда Алексеевна): КР-1 (пусто) -> 0
$ws.Range("H22").Value = 0

# New bold-formatted (empty) cell below the table
$ws.Range("C33").Font.Bold = $true

# Leave the selection where the author's last action left it
$ws.Range("A27:H33").Select()
